$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-6 (existing data) get duplicated as rows 7-11 (identical content,
# same styling). Use targeted PasteSpecial calls so the engine reuses the
# existing cellXfs / numFmt / sharedStrings entries instead of minting new
# ones (which a plain .Copy()/.PasteSpecial() "paste everything" or a naive
# .Value = ... round-trip through COM's numeric auto-coercion would do for
# the numeric-looking text in columns G/H).

for ($i = 0; $i -le 4; $i++) {
    $srcRow = 2 + $i
    $dstRow = 7 + $i

    # Column A: blank cell that only carries a style (s="1") on source rows 2 & 3.
    if ($srcRow -eq 2 -or $srcRow -eq 3) {
        $ws.Cells.Item($srcRow, 1).Copy()
        $ws.Cells.Item($dstRow, 1).PasteSpecial(-4122)   # xlPasteFormats
    }

    # Column B: date serial value styled with the custom date format (s="3").
    # Copy the format (reuses existing xf) then the literal numeric value.
    $ws.Cells.Item($srcRow, 2).Copy()
    $ws.Cells.Item($dstRow, 2).PasteSpecial(-4122)       # xlPasteFormats
    $ws.Cells.Item($srcRow, 2).Copy()
    $ws.Cells.Item($dstRow, 2).PasteSpecial(-4163)       # xlPasteValues

    # Columns C-F: plain shared-string text, default (no) style.
    for ($col = 3; $col -le 6; $col++) {
        $ws.Cells.Item($dstRow, $col).Value = $ws.Cells.Item($srcRow, $col).Value()
    }

    # Columns G-H: shared-string text that LOOKS numeric ("33", "1.77", ...).
    # xlPasteValues preserves the text type instead of coercing to a number.
    for ($col = 7; $col -le 8; $col++) {
        $ws.Cells.Item($srcRow, $col).Copy()
        $ws.Cells.Item($dstRow, $col).PasteSpecial(-4163)  # xlPasteValues
    }
}

$excel.CutCopyMode = 0

$ws.Range("C12").Select()
